# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Estados Unidos (row 9) ---
$ws.Range("B9").Value = 11699
$ws.Range("C9").Value = 2440
$ws.Range("E9").Value = 11417

# --- Noruega (row 17) ---
$ws.Range("B17").Value = 1784
$ws.Range("C17").Value = 193
$ws.Range("E17").Value = 1776

# --- Reunion moves up in the list (new entry ahead of Uzbekistan), pushing
# Uzbekistan / Martinica / Afganistan / Ucrania each one row down; the old
# Reunion row (108) now holds what used to be Ucrania's figures. Update the
# country names and figures for rows 104-108 accordingly.

# Row 104: now Reunion (new data)
$ws.Range("A104").Value = "Reunion"
$ws.Range("B104").Value = 28
$ws.Range("C104").Value = 14
$ws.Range("D104").Value = 0
$ws.Range("E104").Value = 28
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 0

# Row 105: now Uzbekistan
$ws.Range("A105").Value = "Uzbekistan"
$ws.Range("B105").Value = 23
$ws.Range("C105").Value = 5
$ws.Range("D105").Value = 0
$ws.Range("E105").Value = 23
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 0

# Row 106: now Martinica
$ws.Range("A106").Value = "Martinica"
$ws.Range("B106").Value = 23
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 0
$ws.Range("E106").Value = 22
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 1

# Row 107: now Afganistan
$ws.Range("A107").Value = "Afganistan"
$ws.Range("B107").Value = 22
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 1
$ws.Range("E107").Value = 21
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 0

# Row 108: now Ucrania
$ws.Range("A108").Value = "Ucrania"
$ws.Range("B108").Value = 21
$ws.Range("C108").Value = 5
$ws.Range("D108").Value = 0
$ws.Range("E108").Value = 18
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 3

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Marzo de 2020 a las 21:14"
